$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update header label to reflect new "through" date
$ws.Name = "Through 2022-11-11"
$ws.Range("I1").Value = "2022 (through 11-11)"

# Update data values per diff
$ws.Range("I12").Value = 31
$ws.Range("I14").Value = 1430
